# "Generate Report for Handback"
#
# The 7da45bdc-6459-4009-99d2-d152d66656ac source file has now been handed
# back (its translations are in sync with en-US). This updates the status
# on every sheet, records the new "Latest Target File" / "Latest Handback
# File" links and the new "Latest Handback DateTime" for both languages.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # matches styles.xml HyperLink color FF6495ED

function Set-HandbackLink($ws, $cellAddr, $url, $displayText) {
    $ws.Range($cellAddr).Value = $displayText
    $ws.Range($cellAddr).Font.Underline = 2
    $ws.Range($cellAddr).Font.Color = $hyperlinkColor
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $url, "", "", $displayText) | Out-Null
}

# ---------------------------------------------------------------------
# Overview sheet: the 7da45bdc row (row 2) switches from "Ready for
# handoff" to "Handed back: in sync with en-US" for both locale columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"

Set-HandbackLink $wsZhCn "F2" `
    "https://github.com/OpenLocalizationTest/oltest/blob/3ae7c3d900917b165e75055d00c0ec2156e5e522/e2e/7da45bdc-6459-4009-99d2-d152d66656ac.md" `
    "7da45bdc-6459-4009-99d2-d152d66656ac.md"

Set-HandbackLink $wsZhCn "G2" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c6d8243277779c87f8ee055e5bc569d23ce4c29e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7da45bdc-6459-4009-99d2-d152d66656ac.3fef5bedda814bb408c882fca7a3bf4fe37ff90f.zh-cn.xlf" `
    "7da45bdc-6459-4009-99d2-d152d66656ac.3fef5bedda814bb408c882fca7a3bf4fe37ff90f.zh-cn.xlf"

$wsZhCn.Range("H2").Value = "2016-03-24 12:49:32"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

Set-HandbackLink $wsDeDe "F2" `
    "https://github.com/OpenLocalizationTest/oltest/blob/3ae7c3d900917b165e75055d00c0ec2156e5e522/e2e/7da45bdc-6459-4009-99d2-d152d66656ac.md" `
    "7da45bdc-6459-4009-99d2-d152d66656ac.md"

Set-HandbackLink $wsDeDe "G2" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/311c3dbc89a197aaafb86fb38f2fbdb40e28e9f5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7da45bdc-6459-4009-99d2-d152d66656ac.3fef5bedda814bb408c882fca7a3bf4fe37ff90f.de-de.xlf" `
    "7da45bdc-6459-4009-99d2-d152d66656ac.3fef5bedda814bb408c882fca7a3bf4fe37ff90f.de-de.xlf"

$wsDeDe.Range("H2").Value = "2016-03-24 12:49:40"

$wb.Save()
